# console simulator & input screen
$wb = $excel.ActiveWorkbook

$wsTypography = $wb.Worksheets.Item("Typography")
$wsTranslation = $wb.Worksheets.Item("Translation")

# --- Typography sheet: font size for the "Display"/PercentageBarHeadline row changes 28 -> 16 ---
$wsTypography.Range("D20").Value = 16

# --- Translation sheet: split the "Sleep after <value> of inactivity" row into two rows ---
# Insert a new row below the existing SleepWheel_readout row (it inherits formatting from the row above)
$wsTranslation.Rows(67).Insert()

# Row 66 becomes the new "headline" row (was the combined sentence)
$wsTranslation.Range("B66").Value = "SleepWheel_headline"
$wsTranslation.Range("E66").Value = "Sleep"

# Row 67 (newly inserted) becomes the "readout" row, a copy of what row 66 used to contain
$wsTranslation.Range("B67").Value = "SleepWheel_readout"
$wsTranslation.Range("C67").Value = "PercentageBarHeadline"
$wsTranslation.Range("D67").Value = "CENTER"
$wsTranslation.Range("E67").Value = "after <value> of inactivity"

# Extend Table8 (bound to the Translation sheet) by one row to match the new dimension
$lo = $wsTranslation.ListObjects.Item("Table8")
$lo.Resize($wsTranslation.Range("B3:I795"))

# --- View state: Translation tab is now the active tab/sheet, with a new selection ---
$wsTypography.Range("B20").Select()
$wsTranslation.Activate()
$wsTranslation.Range("B61").Select()
